$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "43.185.20"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "2.354.44"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "2.722.97"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "2.348.01"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.797"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "43.168.23"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E28").Value = "  +14.52%  "
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0724"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.46%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -33.00%  "
$ws.Range("D43").Value = "1.943.66"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("E45").Value = "  +3.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "2.586.41"
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.59%  "
